$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.455.24'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.52%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.573.03'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.45%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.25%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.496'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.39%  '

# Row 7
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.07'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.09%  '

# Row 9
$ws.Range("E9").Value = '  -1.22%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0589'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.54%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0865'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.14%  '

# Row 12
$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.801.31'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.14%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.575.79'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.17%  '

# Row 14
$ws.Range("E14").Value = '  -1.12%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.522'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.66%  '

# Row 16
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.06'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.60%  '

# Row 17
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.467.10'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.38%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '213.28'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.02%  '

# Row 19
$ws.Range("D19").Value = '0.0₃0689'
$ws.Range("E19").Value = '  -0.77%  '

# Row 20
$ws.Range("E20").Value = '  -1.22%  '

# Row 21
$ws.Range("E21").Value = '  -0.05%  '

# Row 23
$ws.Range("E23").Value = '  +1.06%  '

# Row 24
$ws.Range("E24").Value = '  +0.52%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.16%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.83'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.03%  '

# Row 27
$ws.Range("E27").Value = '  +0.01%  '

# Row 28
$ws.Range("E28").Value = '  -0.38%  '

# Row 29
$ws.Range("E29").Value = '  -1.64%  '

# Row 30
$ws.Range("E30").Value = '  -0.69%  '

# Row 31
$ws.Range("E31").Value = '  +0.87%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.21'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.24%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.359.90'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.67%  '

# Row 34
$ws.Range("E34").Value = '  -0.58%  '

# Row 35
$ws.Range("E35").Value = '  +1.00%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.971'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.07%  '

# Row 37
$ws.Range("E37").Value = '  +0.21%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0167'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.12%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.532'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.51%  '

# Row 40
$ws.Range("E40").Value = '  +1.41%  '

# Row 41
$ws.Range("E41").Value = '  +0.00%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.974'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.11%  '

# Row 43
$ws.Range("E43").Value = '  -0.38%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '64.28'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.32%  '

# Row 45
$ws.Range("E45").Value = '  -1.13%  '

# Row 46
$ws.Range("E46").Value = '  +0.22%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.711.67'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.20%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.33%  '

# Row 49
$ws.Range("D49").Value = '0.0₇0997'
$ws.Range("E49").Value = '  -0.76%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0956'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.43%  '

# Row 51
$ws.Range("E51").Value = '  -0.41%  '
